# Add descriptions titles 0e2ce4410dd514197a5cd5e0ba72db155e62db5f
#
# - Metadata sheet: fill in the "Title" value (row 5) and "Description" value
#   (row 12), and bump the "Date" value (row 8) to the new generation timestamp.
# - Elements sheet: the root Extension row (row 2) gets its generic
#   "Short"/"Definition" placeholder text replaced by the real title/description,
#   and its "Mapping: RIM Mapping" placeholder ("N/A") is cleared.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B5").Value  = "DMI Classe Risque"
$metadata.Range("B8").Value  = "2026-02-25T08:15:31+00:00"
$metadata.Range("B12").Value = "Extension créée dans ce volet pour représenter la classe de risque."

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("L2").Value  = "DMI Classe Risque"
$elements.Range("M2").Value  = "Extension créée dans ce volet pour représenter la classe de risque."
$elements.Range("AK2").Value = ""
